$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Description text reused from the JD_001 posting (row 2, column C)
$desc = "We are seeking a Software Engineer to build and maintain high-quality software solutions.`nWork with global teams to drive innovation and deliver scalable applications.`nJoin Akkodis and be part of a tech-driven, collaborative environment."

# Add new Job Posting row (row 9) with Job_Id = JD_008
$ws.Cells.Item(9, 1).Value = "JD_008"
$ws.Cells.Item(9, 2).Value = "Senior Analyst "
$ws.Cells.Item(9, 3).Value = $desc
$ws.Cells.Item(9, 4).Value = 1
$ws.Cells.Item(9, 5).Value = 2

# Reset row height auto-fit so no explicit custom row height is persisted
$ws.Cells.Item(9, 3).EntireRow.AutoFit()
